# Removed Test Case Inter-Dependency:
# Rename the product (and its short name) so this loan-product test case no
# longer collides with / depends on another test case's product name, and
# make ProductLoanInput (instead of ProductLoanOutput) the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4227-RBI-EI-DB-SAR-REC-CTRFD-RNI-FEE-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ON-PE-1st"
$newShortName = "422y"

# Update product name (B1) on both the input and output sheets.
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update the short name (B2) on the input sheet - now a text value.
$wsInput.Range("B2").Value = $newShortName

# Make ProductLoanInput the active sheet/tab, with B4 selected.
$wsInput.Activate()
$wsInput.Range("B4").Select()
